$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    @('SingleUseId4064', 'Verdana25', 'Left', 'LTR', '<value>', '<value>', '<value>', '<value>'),
    @('SingleUseId4065', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4066', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4067', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4068', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4069', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4070', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4071', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4072', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4073', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4074', 'Default', 'Center', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4075', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4076', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4077', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4078', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4079', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4080', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4081', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4082', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4083', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4084', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4085', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4086', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4087', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4088', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4089', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4090', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4091', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4092', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4093', 'Default', 'Left', 'LTR', ' ', ' ', ' ', ' '),
    @('SingleUseId4094', 'Small', 'Center', 'LTR', 'Ok', 'Ok', 'Ok', 'Ok'),
    @('SingleUseId4095', 'Small', 'Center', 'LTR', 'Cancel', 'Cancel', 'Cancel', 'Cancel'),
    @('SingleUseId4097', 'Small', 'Center', 'LTR', '?123', '?123', '?123', '?123')
)

$startRow = 656
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $vals = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 6).Value = $vals[4]
    $ws.Cells.Item($r, 7).Value = $vals[5]
    $ws.Cells.Item($r, 8).Value = $vals[6]
    $ws.Cells.Item($r, 9).Value = $vals[7]
}
